$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.475854635238647
$ws.Range("B1").Value = 3.371496438980103
$ws.Range("C1").Value = 4.154712677001953
$ws.Range("D1").Value = 2.226433753967285
$ws.Range("E1").Value = 0.7278503775596619
